# Re-sort TFT experiment results by loss (ascending) and append two new
# in-flight runs (rows 28-29) whose loss is not yet available.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 28,9
$data[0,0] = 'model_month_113_72'
$data[0,1] = 1.618298530578613
$data[0,2] = 30
$data[0,3] = 0.01
$data[0,4] = 36
$data[0,5] = 0.3
$data[0,6] = 22
$data[0,7] = 4
$data[0,8] = 0.01
$data[1,0] = 'model_month_133_72'
$data[1,1] = 2.063173770904541
$data[1,2] = 25
$data[1,3] = 0.01
$data[1,4] = 36
$data[1,5] = 0.3
$data[1,6] = 22
$data[1,7] = 4
$data[1,8] = 0.001
$data[2,0] = 'model_month_365_360'
$data[2,1] = 3.133453845977783
$data[2,2] = 20
$data[2,3] = 0.03
$data[2,4] = 24
$data[2,5] = 0.25
$data[2,6] = 24
$data[2,7] = 4
$data[2,8] = 0.005
$data[3,0] = 'model_month_365_360'
$data[3,1] = 3.346664428710938
$data[3,2] = 20
$data[3,3] = 0.03
$data[3,4] = 24
$data[3,5] = 0.25
$data[3,6] = 24
$data[3,7] = 4
$data[3,8] = 0.01
$data[4,0] = 'model_month_365_168'
$data[4,1] = 3.64
$data[4,2] = 30
$data[4,3] = 0.3118237184625287
$data[4,4] = 13
$data[4,5] = 0.2942142845880926
$data[4,6] = 8
$data[4,7] = 1
$data[4,8] = 0.002283777543696124
$data[5,0] = 'model_month_365_72'
$data[5,1] = 3.702678680419922
$data[5,2] = 25
$data[5,3] = 0.01
$data[5,4] = 36
$data[5,5] = 0.3
$data[5,6] = 22
$data[5,7] = 4
$data[5,8] = 0.001
$data[6,0] = 'model_month_365_168'
$data[6,1] = 3.740181922912598
$data[6,2] = 25
$data[6,3] = 0.01
$data[6,4] = 36
$data[6,5] = 0.3
$data[6,6] = 22
$data[6,7] = 4
$data[6,8] = 0.001
$data[7,0] = 'model_month_200_72'
$data[7,1] = 3.747611522674561
$data[7,2] = 40
$data[7,3] = 0.03
$data[7,4] = 24
$data[7,5] = 0.25
$data[7,6] = 24
$data[7,7] = 2
$data[7,8] = 0.01
$data[8,0] = 'model_month_365_360'
$data[8,1] = 3.845516681671143
$data[8,2] = 40
$data[8,3] = 0.03
$data[8,4] = 24
$data[8,5] = 0.25
$data[8,6] = 24
$data[8,7] = 2
$data[8,8] = 0.01
$data[9,0] = 'model_month_365_168'
$data[9,1] = 4.07178783416748
$data[9,2] = 20
$data[9,3] = 0.05808054277495504
$data[9,4] = 32
$data[9,5] = 0.2997107370422699
$data[9,6] = 19
$data[9,7] = 2
$data[9,8] = 0.15
$data[10,0] = 'model_month_365_72'
$data[10,1] = 4.297078609466553
$data[10,2] = 30
$data[10,3] = 0.03
$data[10,4] = 24
$data[10,5] = 0.25
$data[10,6] = 24
$data[10,7] = 4
$data[10,8] = 0.05
$data[11,0] = 'model_month_365_72'
$data[11,1] = 4.353591442108154
$data[11,2] = 20
$data[11,3] = 0.01
$data[11,4] = 36
$data[11,5] = 0.3
$data[11,6] = 22
$data[11,7] = 4
$data[11,8] = 0.001
$data[12,0] = 'model_month_365_360'
$data[12,1] = 4.596550941467285
$data[12,2] = 20
$data[12,3] = 0.03
$data[12,4] = 24
$data[12,5] = 0.25
$data[12,6] = 24
$data[12,7] = 2
$data[12,8] = 0.1
$data[13,0] = 'model_week_133_72'
$data[13,1] = 5.055328369140625
$data[13,2] = 50
$data[13,3] = 0.01
$data[13,4] = 36
$data[13,5] = 0.3
$data[13,6] = 22
$data[13,7] = 4
$data[13,8] = 0.001
$data[14,0] = 'model_month_365_360'
$data[14,1] = 5.236636638641357
$data[14,2] = 20
$data[14,3] = 0.05808054277495504
$data[14,4] = 32
$data[14,5] = 0.2997107370422699
$data[14,6] = 19
$data[14,7] = 2
$data[14,8] = 0.15
$data[15,0] = 'model_month_700_360'
$data[15,1] = 5.345937252044678
$data[15,2] = 30
$data[15,3] = 0.05808054277495504
$data[15,4] = 32
$data[15,5] = 0.2997107370422699
$data[15,6] = 19
$data[15,7] = 2
$data[15,8] = 0.01
$data[16,0] = 'model_month_700_504'
$data[16,1] = 5.529139041900635
$data[16,2] = 30
$data[16,3] = 0.05808054277495504
$data[16,4] = 32
$data[16,5] = 0.2997107370422699
$data[16,6] = 19
$data[16,7] = 2
$data[16,8] = 0.01
$data[17,0] = 'model_month_700_168'
$data[17,1] = 5.57568883895874
$data[17,2] = 20
$data[17,3] = 0.05808054277495504
$data[17,4] = 32
$data[17,5] = 0.2997107370422699
$data[17,6] = 19
$data[17,7] = 2
$data[17,8] = 0.05
$data[18,0] = 'model_week_133_24'
$data[18,1] = 6.060751438140869
$data[18,2] = 100
$data[18,3] = 0.01
$data[18,4] = 36
$data[18,5] = 0.3
$data[18,6] = 22
$data[18,7] = 4
$data[18,8] = 0.001
$data[19,0] = 'model_month_700_504'
$data[19,1] = 6.403566360473633
$data[19,2] = 25
$data[19,3] = 0.05
$data[19,4] = 32
$data[19,5] = 0.25
$data[19,6] = 32
$data[19,7] = 6
$data[19,8] = 0.01
$data[20,0] = 'model_week_133_72'
$data[20,1] = 8.389998435974121
$data[20,2] = 25
$data[20,3] = 0.01
$data[20,4] = 36
$data[20,5] = 0.3
$data[20,6] = 22
$data[20,7] = 4
$data[20,8] = 0.001
$data[21,0] = 'model_week_365_144'
$data[21,1] = 8.747285842895508
$data[21,2] = 30
$data[21,3] = 0.4696487911301886
$data[21,4] = 8
$data[21,5] = 0.1039155309883323
$data[21,6] = 8
$data[21,7] = 4
$data[21,8] = 0.05
$data[22,0] = 'model_week_365_72'
$data[22,1] = 10.83476448059082
$data[22,2] = 30
$data[22,3] = 0.4696487911301886
$data[22,4] = 8
$data[22,5] = 0.1039155309883323
$data[22,6] = 8
$data[22,7] = 4
$data[22,8] = 0.05
$data[23,0] = 'model_week_700_144'
$data[23,1] = 12.59808349609375
$data[23,2] = 30
$data[23,3] = 0.01058413796901732
$data[23,4] = 36
$data[23,5] = 0.2960659121982872
$data[23,6] = 22
$data[23,7] = 4
$data[23,8] = 0.05
$data[24,0] = 'model_week_700_72'
$data[24,1] = 12.97217750549316
$data[24,2] = 30
$data[24,3] = 0.01058413796901732
$data[24,4] = 36
$data[24,5] = 0.2960659121982872
$data[24,6] = 22
$data[24,7] = 4
$data[24,8] = 0.05
$data[25,0] = 'model_week_133_24'
$data[25,1] = 13.58229351043701
$data[25,2] = 25
$data[25,3] = 0.01
$data[25,4] = 36
$data[25,5] = 0.3
$data[25,6] = 22
$data[25,7] = 4
$data[25,8] = 0.001
$data[26,0] = 'model_month_365_360'
$data[26,1] = $null
$data[26,2] = 30
$data[26,3] = 0.03
$data[26,4] = 24
$data[26,5] = 0.25
$data[26,6] = 24
$data[26,7] = 2
$data[26,8] = 0.01
$data[27,0] = 'model_month_365_360'
$data[27,1] = $null
$data[27,2] = 20
$data[27,3] = 0.03
$data[27,4] = 24
$data[27,5] = 0.25
$data[27,6] = 24
$data[27,7] = 2
$data[27,8] = 0.1

$ws.Range("A2:I29").Value = $data
